$d = $word.ActiveDocument

$replacements = @(
    @("71×54=", "37×21="),
    @("77×46=", "71×95="),
    @("76×88=", "98×98="),
    @("96×38=", "89×25="),
    @("30×64=", "93×34="),
    @("58×31=", "32×76="),
    @("66×32=", "27×20="),
    @("12×88=", "57×74="),
    @("54×91=", "91×12="),
    @("98×75=", "61×19="),
    @("25×17=", "60×62="),
    @("63×99=", "88×46="),
    @("64×60=", "36×78="),
    @("92×54=", "73×87="),
    @("23×87=", "71×94="),
    @("35×44=", "30×31="),
    @("35×90=", "49×86="),
    @("28×61=", "93×18="),
    @("85×26=", "89×34="),
    @("41×97=", "93×65="),
    @("58×69=", "56×67="),
    @("39×35=", "51×98="),
    @("77×81=", "35×91="),
    @("39×29=", "75×59="),
    @("96×46=", "13×68=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
